# Adds a new column Z ("11-10-2020") to the COVID19_TIMESERIESDATA sheet,
# mirroring the formatting of the existing last data column (Y), and fills
# in the per-state active-case counts for that date (rows 2-36).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell Z1: matches the style of Y1 (bold, centered, bordered) ---
# Copy Y1's formatting onto Z1 first so Z1 looks the same as the rest of the
# date header row, then write the value as literal text. The leading
# apostrophe forces a text entry so "11-10-2020" is stored as a string
# rather than being auto-parsed into a date serial.
$ws.Range("Y1").Copy()
$ws.Range("Z1").PasteSpecial(-4122)
$ws.Range("Z1").Value = "'11-10-2020"

# --- Data cells Z2:Z36: plain numeric active-case counts for 11-10-2020 ---
$zValues = @{
    2  = 193
    3  = 46624
    4  = 2940
    5  = 29221
    6  = 11165
    7  = 1229
    8  = 27369
    9  = 100
    10 = 22007
    11 = 4658
    12 = 15936
    13 = 10677
    14 = 2718
    15 = 10796
    16 = 8362
    17 = 120948
    18 = 96003
    19 = 1022
    20 = 15612
    21 = 221615
    22 = 2608
    23 = 2437
    24 = 191
    25 = 1238
    26 = 24414
    27 = 4719
    28 = 9752
    29 = 21354
    30 = 450
    31 = 44150
    32 = 25713
    33 = 3951
    34 = 7321
    35 = 40210
    36 = 29793
}

foreach ($row in $zValues.Keys) {
    $ws.Cells.Item($row, 26).Value = $zValues[$row]
}
